$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new person "Jan Student" in row 6 (continuing the list started in A3:A5)
$ws.Range("B6").Value = "Jan"
$ws.Range("D6").Value = "Student"

# Re-establish the shared formula so it spans A6:A7 (mirrors dragging the
# fill handle of A5 down through A7), matching rows 3-5's original pattern.
$ws.Range("A6:A7").Formula = '=IF($B6="","",CONCATENATE($B6," ",IF($C6="",D6,CONCATENATE($C6," ",$D6))))'

# Update the selected cell to match the author's final cursor position
$ws.Range("E5").Select()
